$p = $ppt.ActivePresentation
$sm = $p.SlideMaster
$theme = $sm.Theme
$tcs = $theme.ThemeColorScheme
try {
    $tcs.BogusXYZ = "hello"
    Write-Host "set bogus OK"
} catch {
    Write-Host ("ERROR: " + $_.Exception.Message)
}
